# Updated cryptos list on Sun Nov  5 02:37:39 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to be treated as plain text so numeric-looking
    # strings (e.g. "242.55") are not silently coerced into floating
    # point numbers (e.g. 242.55000000000001) by Excel's auto-detection.
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
}

# Row 13 / 14: Chainlink and WrappedEther swap positions (row 13 becomes Chainlink, row 14 becomes WrappedEther)
Set-TextValue "B13" "Chainlink"
Set-TextValue "C13" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D13" "11.84"
Set-TextValue "E13" "  +3.46%  "

Set-TextValue "B14" "WrappedEther"
Set-TextValue "C14" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D14" "1.881.58"
Set-TextValue "E14" "  +1.50%  "

# Row 33 / 34: ImmutableX and InternetComputer(DFINITY) swap positions
Set-TextValue "B33" "ImmutableX"
Set-TextValue "C33" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D33" "0.929"
Set-TextValue "E33" "  +21.63%  "

Set-TextValue "B34" "InternetComputer(DFINITY)"
Set-TextValue "C34" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D34" "4.08"
Set-TextValue "E34" "  +2.46%  "

# Row 46 / 47: Gas and HuobiToken swap positions
Set-TextValue "B46" "Gas"
Set-TextValue "C46" "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
Set-TextValue "D46" "12.62"
Set-TextValue "E46" "  +44.13%  "

Set-TextValue "B47" "HuobiToken"
Set-TextValue "C47" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D47" "2.43"
Set-TextValue "E47" "  +0.23%  "

# Remaining rows: update Price (D) and/or Volume(1h) (E) values only

Set-TextValue "D2" "35.461.20"
Set-TextValue "E2" "  +1.13%  "

Set-TextValue "D3" "1.873.19"
Set-TextValue "E3" "  +1.08%  "

Set-TextValue "D4" "1.01"
Set-TextValue "E4" "  +0.43%  "

Set-TextValue "D5" "242.55"
Set-TextValue "E5" "  +4.06%  "

Set-TextValue "E6" "  +2.06%  "

Set-TextValue "E7" "  +0.42%  "

Set-TextValue "D8" "43.36"
Set-TextValue "E8" "  +6.21%  "

Set-TextValue "E9" "  +0.01%  "

Set-TextValue "E10" "  +1.26%  "

Set-TextValue "D11" "0.0993"
Set-TextValue "E11" "  +0.90%  "

Set-TextValue "D12" "2.145.69"
Set-TextValue "E12" "  +1.14%  "

Set-TextValue "D15" "0.685"
Set-TextValue "E15" "  +1.24%  "

Set-TextValue "E16" "  +1.85%  "

Set-TextValue "D17" "35.500.92"
Set-TextValue "E17" "  +1.07%  "

Set-TextValue "D18" "71.01"
Set-TextValue "E18" "  +1.15%  "

Set-TextValue "E19" "  +1.37%  "

Set-TextValue "D20" "242.64"
Set-TextValue "E20" "  +0.74%  "

Set-TextValue "D21" "12.33"
Set-TextValue "E21" "  +0.29%  "

Set-TextValue "E22" "  +1.24%  "

Set-TextValue "E23" "  +0.45%  "

Set-TextValue "D24" "2.26"
Set-TextValue "E24" "  +0.30%  "

Set-TextValue "D25" "170.98"
Set-TextValue "E25" "  -1.29%  "

Set-TextValue "D26" "2.04"
Set-TextValue "E26" "  +29.07%  "

Set-TextValue "E27" "  +4.94%  "

Set-TextValue "D28" "17.82"
Set-TextValue "E28" "  +1.41%  "

Set-TextValue "E29" "  +0.67%  "

Set-TextValue "E30" "  +1.37%  "

Set-TextValue "E31" "  +0.63%  "

Set-TextValue "D32" "4.05"
Set-TextValue "E32" "  +2.27%  "

Set-TextValue "D35" "1.78"
Set-TextValue "E35" "  +11.40%  "

Set-TextValue "D36" "2.06"
Set-TextValue "E36" "  +4.29%  "

Set-TextValue "D37" "1.36"
Set-TextValue "E37" "  +11.19%  "

Set-TextValue "E38" "  +1.87%  "

Set-TextValue "E39" "  +3.69%  "

Set-TextValue "D40" "89.57"
Set-TextValue "E40" "  -1.04%  "

Set-TextValue "D41" "1.353.99"
Set-TextValue "E41" "  +0.07%  "

Set-TextValue "D42" "15.24"
Set-TextValue "E42" "  +3.79%  "

Set-TextValue "D43" "0.0591"
Set-TextValue "E43" "  +11.26%  "

Set-TextValue "E44" "  +3.12%  "

Set-TextValue "D45" "47.14"
Set-TextValue "E45" "  +38.55%  "

Set-TextValue "D48" "6.71"
Set-TextValue "E48" "  +5.77%  "

Set-TextValue "E49" "  -1.37%  "

Set-TextValue "D50" "2.062.47"
Set-TextValue "E50" "  +0.94%  "

Set-TextValue "D51" "0.0685"
Set-TextValue "E51" "  +2.40%  "
